$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.121.00'
$ws.Range("E2").Value = '  -1.88%  '

$ws.Range("D3").Value = '2.906.08'
$ws.Range("E3").Value = '  -0.83%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '347.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.549'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.16%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.602'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.29'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.65%  '

$ws.Range("E11").Value = '  +1.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0839'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.70'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.58%  '

$ws.Range("D14").Value = '3.366.65'
$ws.Range("E14").Value = '  -0.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.30%  '

$ws.Range("D16").Value = '2.893.90'
$ws.Range("E16").Value = '  -0.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.949'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = '51.123.63'
$ws.Range("E18").Value = '  -1.94%  '

$ws.Range("E19").Value = '  +3.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.73%  '

$ws.Range("E21").Value = '  -6.49%  '

$ws.Range("D22").Value = '0.0₃0955'
$ws.Range("E22").Value = '  -2.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.44'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '260.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.25%  '

$ws.Range("E25").Value = '  -4.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.170'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.45%  '

$ws.Range("E27").Value = '  +0.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '26.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.36'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.103'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.71%  '

$ws.Range("E32").Value = '  +0.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.10%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.13'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.11'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.36%  '

$ws.Range("E40").Value = '  -6.18%  '

$ws.Range("E41").Value = '  -3.02%  '

$ws.Range("E42").Value = '  -2.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.65'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.57%  '

$ws.Range("E45").Value = '  -2.79%  '

$ws.Range("D46").Value = '2.083.98'
$ws.Range("E46").Value = '  -4.88%  '

$ws.Range("E47").Value = '  -6.90%  '

$ws.Range("E48").Value = '  -10.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.236'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.99%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0331'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.880'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.48%  '
